# Update column G ("K") values for rows 2-17 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
